# Weekly price-list refresh: existing weekly rows (53-78) shift to new
# dates/values, and three brand-new rows (79-81) are appended for the
# 2021-08-17 (Region de OHiggins) entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @{}
$rowsData[53] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44455, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Fuji royal', 'Calibre 80', 270, 20000, 21000, 20500, '$/caja 18 kilos embalada', 'Región de O''Higgins', 1139, 18)
$rowsData[54] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44455, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Granny Smith', 'Calibre 90', 300, 20000, 21000, 20500, '$/caja 18 kilos embalada', 'Región de O''Higgins', 1139, 18)
$rowsData[55] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44455, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Scarlett', 'Calibre 90', 300, 20000, 21000, 20500, '$/caja 18 kilos embalada', 'Región de O''Higgins', 1139, 18)
$rowsData[56] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44280, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Fuji royal', 'Segunda', 300, 19000, 20000, 19500, '$/caja 18 kilos granel', 'Región de O''Higgins', 1083, 18)
$rowsData[57] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44280, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Granny Smith', 'Segunda', 250, 19000, 20000, 19500, '$/caja 18 kilos granel', 'Región de O''Higgins', 1083, 18)
$rowsData[58] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44390, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Granny Smith', 'Calibre 80', 250, 15000, 16000, 15500, '$/caja 18 kilos embalada', 'Provincia de Curicó', 861, 18)
$rowsData[59] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44390, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Pink Lady', 'Calibre 80', 200, 15000, 16000, 15500, '$/caja 18 kilos embalada', 'Provincia de Curicó', 861, 18)
$rowsData[60] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44390, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Royal Gala', 'Calibre 80', 150, 15000, 16000, 15500, '$/caja 18 kilos embalada', 'Provincia de Curicó', 861, 18)
$rowsData[61] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44390, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Royal Gala', 'Calibre 90', 150, 15000, 16000, 15500, '$/caja 18 kilos embalada', 'Provincia de Curicó', 861, 18)
$rowsData[62] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44390, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Scarlett', 'Calibre 90', 300, 15000, 16000, 15500, '$/caja 18 kilos embalada', 'Provincia de Curicó', 861, 18)
$rowsData[63] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44308, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Fuji royal', 'Segunda', 250, 17000, 18000, 17500, '$/caja 18 kilos granel', 'Región de O''Higgins', 972, 18)
$rowsData[64] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44308, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Granny Smith', 'Segunda', 300, 17000, 18000, 17500, '$/caja 18 kilos granel', 'Región de O''Higgins', 972, 18)
$rowsData[65] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44308, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Royal Gala', 'Segunda', 250, 17000, 18000, 17500, '$/caja 18 kilos granel', 'Región de O''Higgins', 972, 18)
$rowsData[66] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44166, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Fuji royal', 'Calibre 80', 250, 24000, 25000, 24500, '$/caja 18 kilos embalada', 'Provincia de Curicó', 1361, 18)
$rowsData[67] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44166, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Granny Smith', 'Calibre 80', 300, 24000, 25000, 24500, '$/caja 18 kilos embalada', 'Provincia de Curicó', 1361, 18)
$rowsData[68] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44397, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Fuji royal', 'Segunda', 270, 16000, 17000, 16500, '$/caja 18 kilos granel', 'Región de O''Higgins', 917, 18)
$rowsData[69] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44397, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Granny Smith', 'Segunda', 300, 16000, 17000, 16500, '$/caja 18 kilos granel', 'Región de O''Higgins', 917, 18)
$rowsData[70] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44397, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Royal Gala', 'Segunda', 250, 16000, 17000, 16500, '$/caja 18 kilos granel', 'Región de O''Higgins', 917, 18)
$rowsData[71] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44397, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Scarlett', 'Segunda', 300, 16000, 17000, 16500, '$/caja 18 kilos granel', 'Región de O''Higgins', 917, 18)
$rowsData[72] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44351, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Fuji royal', 'Segunda', 300, 16000, 17000, 16500, '$/caja 18 kilos granel', 'Provincia de Curicó', 917, 18)
$rowsData[73] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44351, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Granny Smith', 'Segunda', 300, 16000, 17000, 16500, '$/caja 18 kilos granel', 'Provincia de Curicó', 917, 18)
$rowsData[74] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44351, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Royal Gala', 'Segunda', 250, 16000, 17000, 16500, '$/caja 18 kilos granel', 'Provincia de Curicó', 917, 18)
$rowsData[75] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44411, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Fuji royal', 'Segunda', 250, 16000, 17000, 16500, '$/caja 18 kilos granel', 'Región de O''Higgins', 917, 18)
$rowsData[76] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44411, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Granny Smith', 'Segunda', 300, 16000, 17000, 16500, '$/caja 18 kilos granel', 'Región de O''Higgins', 917, 18)
$rowsData[77] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44411, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Royal Gala', 'Segunda', 270, 16000, 17000, 16500, '$/caja 18 kilos granel', 'Región de O''Higgins', 917, 18)
$rowsData[78] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44411, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Scarlett', 'Segunda', 250, 16000, 17000, 16500, '$/caja 18 kilos granel', 'Región de O''Higgins', 917, 18)
$rowsData[79] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44425, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Granny Smith', 'Segunda', 300, 17000, 18000, 17500, '$/caja 18 kilos granel', 'Región de O''Higgins', 972, 18)
$rowsData[80] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44425, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Pink Lady', 'Segunda', 270, 17000, 18000, 17500, '$/caja 18 kilos granel', 'Región de O''Higgins', 972, 18)
$rowsData[81] = @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44425, 15, 'Fruta', 100104, 'Frutos de pepita', 100104002, 'Manzana', 'Royal Gala', 'Segunda', 300, 17000, 18000, 17500, '$/caja 18 kilos granel', 'Región de O''Higgins', 972, 18)

# Rows that are brand-new (not present before) need the date column's
# number format applied explicitly -- existing rows already carry it.
$newRows = @(79, 80, 81)

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
    if ($newRows -contains $r) {
        $ws.Cells.Item($r, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
    }
}
